$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (2021年) is a brand-new row appended right after the existing last
# row (row 11, 2020年). Copy row 11's column-A formatting (bold font, thin
# border, center/top alignment - style index 1 in the original file) onto
# A12 before writing its value, so the new label cell matches the look of
# every other year label in column A.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 13
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 1973
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = 54
$ws.Range("I12").Value = 8
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 6
$ws.Range("L12").Value = 15
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = ""
$ws.Range("P12").Value = ""
$ws.Range("Q12").Value = 146
$ws.Range("R12").Value = 1
$ws.Range("S12").Value = ""
$ws.Range("T12").Value = 10
$ws.Range("U12").Value = ""
$ws.Range("V12").Value = ""
$ws.Range("W12").Value = 87
$ws.Range("X12").Value = ""
$ws.Range("Y12").Value = 1622
